$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("BDSBaPCF")

# Set biomass (row 9) and municipal solid waste (row 17) plant types to 0
$ws.Range("B9").Value = 0
$ws.Range("B17").Value = 0

# Move the BDSBaPCF sheet's remembered selection from B14 to B18, then
# restore "About" as the active sheet (matches original active tab).
$ws.Activate()
$ws.Range("B18").Select()
$wsAbout.Activate()
